$wb = $excel.ActiveWorkbook

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4053.6
$ws.Range("J32").Value = 4567
$ws.Range("L32").Value = 4567
$ws.Range("N32").Value = -5219

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 29072.77
$ws.Range("I131").Value = 2060.348
$ws.Range("J131").Value = 236168
$ws.Range("K131").Value = 6181.044
$ws.Range("L131").Value = 708504
$ws.Range("M131").Value = -1141.044
$ws.Range("N131").Value = -718584

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6525.625
$ws.Range("I137").Value = 10606.917
$ws.Range("J137").Value = 2444.3333
$ws.Range("K137").Value = 31820.751
$ws.Range("L137").Value = 7332.999899999999
$ws.Range("M137").Value = -29270.751
$ws.Range("N137").Value = -12432.9999

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2704.47
$ws.Range("I138").Value = 2165.4482
$ws.Range("J138").Value = 2924.6338
$ws.Range("K138").Value = 6496.344599999999
$ws.Range("L138").Value = 8773.901400000001
$ws.Range("M138").Value = -1356.344599999999
$ws.Range("N138").Value = -19053.9014

# ALC row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 74993
$ws.Range("J139").Value = 74993
$ws.Range("L139").Value = 74993
$ws.Range("N139").Value = -85273

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1575.4762
$ws.Range("I45").Value = 1037.0834
$ws.Range("K45").Value = 1037.0834
$ws.Range("M45").Value = -660.0834

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3101.8965
$ws.Range("I61").Value = 1889.1428
$ws.Range("K61").Value = 1889.1428
$ws.Range("M61").Value = -1677.1428

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3101.8965
$ws.Range("I136").Value = 1889.1428
$ws.Range("K136").Value = 5667.428400000001
$ws.Range("M136").Value = -3117.428400000001

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 74999
$ws.Range("J139").Value = 74999
$ws.Range("L139").Value = 74999
$ws.Range("N139").Value = -85279

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1891
$ws.Range("I16").Value = 1733.4
$ws.Range("K16").Value = 1733.4
$ws.Range("M16").Value = -1446.4

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25022268
$ws.Range("J31").Value = 62552976
$ws.Range("L31").Value = 62552976
$ws.Range("N31").Value = -62553566

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 25022268
$ws.Range("J34").Value = 62552976
$ws.Range("L34").Value = 62552976
$ws.Range("N34").Value = -62553380

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1608.3125
$ws.Range("I94").Value = 1315.7142
$ws.Range("J94").Value = 1835.8889
$ws.Range("K94").Value = 1315.7142
$ws.Range("L94").Value = 1835.8889
$ws.Range("M94").Value = -864.7141999999999
$ws.Range("N94").Value = -2737.8889

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4222.7207
$ws.Range("I105").Value = 1671.4286
$ws.Range("J105").Value = 6658.0454
$ws.Range("K105").Value = 1671.4286
$ws.Range("L105").Value = 6658.0454
$ws.Range("M105").Value = 75.57140000000004
$ws.Range("N105").Value = -10152.0454

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1891
$ws.Range("I113").Value = 1733.4
$ws.Range("K113").Value = 1733.4
$ws.Range("M113").Value = 436.5999999999999

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3240.4167
$ws.Range("I132").Value = 2766.55
$ws.Range("J132").Value = 5609.75
$ws.Range("K132").Value = 8299.650000000001
$ws.Range("L132").Value = 16829.25
$ws.Range("M132").Value = -5769.650000000001
$ws.Range("N132").Value = -21889.25

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4912.6333
$ws.Range("I134").Value = 4944.1035
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 14832.3105
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -12297.3105
$ws.Range("N134").Value = -17070

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 20834046
$ws.Range("I68").Value = 27778312
$ws.Range("K68").Value = 83334936
$ws.Range("M68").Value = -83334125

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 20834046
$ws.Range("I71").Value = 27778312
$ws.Range("K71").Value = 250004808
$ws.Range("M71").Value = -250000752

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 854.1
$ws.Range("J122").Value = 989.6
$ws.Range("L122").Value = 8906.4
$ws.Range("N122").Value = -13806.4

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7938171.5
$ws.Range("J131").Value = 1776.3928
$ws.Range("L131").Value = 5329.178400000001
$ws.Range("N131").Value = -15409.1784

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3431
$ws.Range("I141").Value = 3403.4
$ws.Range("K141").Value = 10210.2
$ws.Range("M141").Value = -5030.200000000001

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 10135
$ws.Range("I126").Value = 12136.333
$ws.Range("J126").Value = 8848.429
$ws.Range("K126").Value = 36408.999
$ws.Range("L126").Value = 26545.287
$ws.Range("M126").Value = -33938.999
$ws.Range("N126").Value = -31485.287

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 925348.9399999999
$ws.Range("I40").Value = 1016908.6
$ws.Range("J40").Value = 9751.5
$ws.Range("K40").Value = 1016908.6
$ws.Range("L40").Value = 9751.5
$ws.Range("M40").Value = -1016772.6
$ws.Range("N40").Value = -10023.5

# LTW row 87
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# LTW row 90
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2697.3103
$ws.Range("I122").Value = 2574.5652
$ws.Range("J122").Value = 3167.8333
$ws.Range("K122").Value = 7723.6956
$ws.Range("L122").Value = 9503.499899999999
$ws.Range("M122").Value = -5273.6956
$ws.Range("N122").Value = -14403.4999

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 651
$ws.Range("I107").Value = 607.8
$ws.Range("J107").Value = 687
$ws.Range("K107").Value = 1823.4
$ws.Range("L107").Value = 2061
$ws.Range("M107").Value = 96.60000000000014
$ws.Range("N107").Value = -5901

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 978.56665
$ws.Range("I113").Value = 847.9
$ws.Range("K113").Value = 2543.7
$ws.Range("M113").Value = -373.6999999999998

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1600.2903
$ws.Range("I122").Value = 1340.3043
$ws.Range("K122").Value = 4020.9129
$ws.Range("M122").Value = -1570.9129
